$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replaces paragraph $index's content with raw WordprocessingML, after checking
# its current plain text matches $expectText (sanity guard against drift).
function Set-ParaXml($index, $expectText, $inner) {
    $para = $d.Paragraphs($index)
    $actual = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($actual -ne $expectText) {
        Write-Output "WARNING paragraph $index text mismatch: expected [$expectText] got [$actual]"
    }
    [void]$para.Range.InsertXML("<w:p $wns>$inner</w:p>")
}

# 1. "Lighthouse" -> mark as a spell-check exception (proofErr spellStart/spellEnd)
Set-ParaXml 1 "Lighthouse" @"
<w:proofErr w:type="spellStart"/><w:r><w:t>Lighthouse</w:t></w:r><w:proofErr w:type="spellEnd"/>
"@

# 3. "Passer aux formats de nouvelle génération (WebP)" -> isolate "WebP" with proofErr
Set-ParaXml 3 "Passer aux formats de nouvelle génération (WebP)" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Passer aux formats de nouvelle génération (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WebP</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r>
"@

# 6. "Différer le chargement d'image hors écran" -> bold
Set-ParaXml 6 "Différer le chargement d’image hors écran" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Différer le chargement d’image hors écran</w:t></w:r>
"@

# 7. "Réduire les ressources CSS inutilisées" -> bold
Set-ParaXml 7 "Réduire les ressources CSS inutilisées" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Réduire les ressources CSS inutilisées</w:t></w:r>
"@

# 8. "Réduire les ressources JS inutilisées" -> bold
Set-ParaXml 8 "Réduire les ressources JS inutilisées" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Réduire les ressources JS inutilisées</w:t></w:r>
"@

# 10. "Wave" -> mark as a spell-check exception
Set-ParaXml 10 "Wave" @"
<w:proofErr w:type="spellStart"/><w:r><w:t>Wave</w:t></w:r><w:proofErr w:type="spellEnd"/>
"@

# 12. "Ajouter les label aux balises form" -> split runs; flag the grammar ("les label")
#     and the spelling exception ("form")
Set-ParaXml 12 "Ajouter les label aux balises form" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ajouter </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>les label</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> aux balises </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>form</w:t></w:r><w:proofErr w:type="spellEnd"/>
"@

# 13. "Modifier le title" -> split runs; flag the spelling exception ("title")
Set-ParaXml 13 "Modifier le title" @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Modifier le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>title</w:t></w:r><w:proofErr w:type="spellEnd"/>
"@

# 16. "Vérifier les éléments de région (header, nav, footer…) de la page." -> split runs;
#     flag the spelling exceptions ("nav" and "footer")
Set-ParaXml 16 "Vérifier les éléments de région (header, nav, footer…) de la page." @"
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Vérifier les éléments de région (header, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>footer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…) de la page.</w:t></w:r>
"@

Write-Output "done"
